$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.225.68'
$ws.Range('E2').Value = '  +2.53%  '

$ws.Range('D3').Value = '1.439.43'
$ws.Range('E3').Value = '  +3.89%  '

$ws.Range('E4').Value = '  +0.37%  '

$ws.Range('B5').Value = 'USDC'
$ws.Range('C5').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D5').Value = '''0.9072'
$ws.Range('E5').Value = '  -9.60%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '''277.39'
$ws.Range('E6').Value = '  +3.30%  '

$ws.Range('D7').Value = '''0.3643'
$ws.Range('E7').Value = '  +0.14%  '

$ws.Range('E8').Value = '  +2.46%  '

$ws.Range('D9').Value = '''39.02'
$ws.Range('E9').Value = '  +2.21%  '

$ws.Range('D10').Value = '''1.018'
$ws.Range('E10').Value = '  +5.03%  '

$ws.Range('D11').Value = '''0.06531'
$ws.Range('E11').Value = '  +2.32%  '

$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  -0.30%  '

$ws.Range('D13').Value = '''5.373'
$ws.Range('E13').Value = '  +2.08%  '

$ws.Range('D14').Value = '''17.56'
$ws.Range('E14').Value = '  +7.52%  '

$ws.Range('D15').Value = '''6.051'
$ws.Range('E15').Value = '  +0.35%  '

$ws.Range('D16').Value = '''0.00001017'
$ws.Range('E16').Value = '  +3.09%  '

$ws.Range('D17').Value = '1.439.79'
$ws.Range('E17').Value = '  +3.74%  '

$ws.Range('D18').Value = '''0.9430'
$ws.Range('E18').Value = '  -6.03%  '

$ws.Range('D19').Value = '''0.05653'
$ws.Range('E19').Value = '  +0.74%  '

$ws.Range('D20').Value = '''67.65'
$ws.Range('E20').Value = '  -2.87%  '

$ws.Range('D21').Value = '''5.388'
$ws.Range('E21').Value = '  -1.75%  '

$ws.Range('D22').Value = '''14.38'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('D23').Value = '''10.77'
$ws.Range('E23').Value = '  +3.02%  '

$ws.Range('D24').Value = '''2.235'
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('D25').Value = '20.281.39'
$ws.Range('E25').Value = '  +2.79%  '

$ws.Range('D26').Value = '''2.157'
$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').Value = '''137.71'
$ws.Range('E27').Value = '  +1.44%  '

$ws.Range('D28').Value = '''16.91'
$ws.Range('E28').Value = '  +2.69%  '

$ws.Range('D29').Value = '1.591.95'
$ws.Range('E29').Value = '  +2.98%  '

$ws.Range('D30').Value = '''109.88'
$ws.Range('E30').Value = '  +2.69%  '

$ws.Range('D31').Value = '''3.887'
$ws.Range('E31').Value = '  +2.03%  '

$ws.Range('D32').Value = '''0.7983'
$ws.Range('E32').Value = '  +2.34%  '

$ws.Range('D33').Value = '''4.800'
$ws.Range('E33').Value = '  -7.45%  '

$ws.Range('D34').Value = '''0.07678'
$ws.Range('E34').Value = '  +1.51%  '

$ws.Range('D35').Value = '''0.05921'
$ws.Range('E35').Value = '  +6.38%  '

$ws.Range('D36').Value = '''1.439'
$ws.Range('E36').Value = '  +11.86%  '

$ws.Range('D37').Value = '''1.146'
$ws.Range('E37').Value = '  +10.47%  '

$ws.Range('D38').Value = '''4.638'
$ws.Range('E38').Value = '  -0.77%  '

$ws.Range('D39').Value = '''0.01984'
$ws.Range('E39').Value = '  -0.97%  '

$ws.Range('D40').Value = '''10.14'
$ws.Range('E40').Value = '  +2.11%  '

$ws.Range('D41').Value = '''0.1835'
$ws.Range('E41').Value = '  -1.39%  '

$ws.Range('D42').Value = '''0.9134'
$ws.Range('E42').Value = '  -8.98%  '

$ws.Range('D43').Value = '''7.069'
$ws.Range('E43').Value = '  -13.47%  '

$ws.Range('D44').Value = '''3.514'
$ws.Range('E44').Value = '  +1.61%  '

$ws.Range('D45').Value = '''0.5225'
$ws.Range('E45').Value = '  +1.72%  '

$ws.Range('D46').Value = '''12.00'
$ws.Range('E46').Value = '  +2.31%  '

$ws.Range('D47').Value = '''118.30'
$ws.Range('E47').Value = '  +9.93%  '

$ws.Range('D48').Value = '''0.5131'
$ws.Range('E48').Value = '  +4.08%  '

$ws.Range('D49').Value = '''1.754'
$ws.Range('E49').Value = '  +2.81%  '

$ws.Range('E50').Value = '  +4.82%  '

$ws.Range('D51').Value = '''0.9886'
$ws.Range('E51').Value = '  -1.59%  '
